# Adapt column header formatting to respective input file names.
#
# The sheet has 21 header columns in row 1:
#   A..J  -> "<Name>_old" columns (describe the "before"/FV2210 input file)
#   K     -> "diff"
#   L..U  -> "<Name>_new" columns (describe the "after"/FV2304 input file)
#
# The headers are renamed so the "_old" suffix becomes "_FV2210" and the
# "_new" suffix becomes "_FV2304". A table (ListObject) is added over the
# whole data range using those (new) header names, and the first row is
# frozen so the header stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (row 1) ----------------------------------

$headersFV2210 = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)

$headersFV2304 = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

# Columns A (1) .. J (10)
for ($i = 0; $i -lt $headersFV2210.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headersFV2210[$i]
}

# Column K (11) is "diff" and stays unchanged.

# Columns L (12) .. U (21)
for ($i = 0; $i -lt $headersFV2304.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value2 = $headersFV2304[$i]
}

# --- 2. Turn the used range into an Excel Table (ListObject) -------------

$dataRange = $ws.Range("A1:U82")
$table = $ws.ListObjects.Add(1, $dataRange, [System.Reflection.Missing]::Value, 1)
$table.Name = "Table1"

# --- 3. Freeze the header row ---------------------------------------------

[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
